# Update cryptocurrency price/volume figures on Sheet1 to reflect the
# refreshed snapshot pulled by the GitHub Actions symbol-list update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value (Price column D, Volume(1h) column E).
# NumberFormat is forced to Text ("@") before writing so that values such as
# "328.94", "0.00002104" or "4.64%" are stored verbatim as text, matching the
# original inline-string cell contents instead of being reinterpreted as numbers.
$updates = [ordered]@{
    "D2" = "328.94"
    "E2" = "4.64%"
    "D3" = "40.50"
    "E3" = "9.07%"
    "D4" = "5.740"
    "E4" = "12.02%"
    "D5" = "0.08113"
    "E5" = "2.39%"
    "D6" = "4.604"
    "E6" = "3.96%"
    "D7" = "8.769"
    "E7" = "4.05%"
    "D8" = "1.969"
    "E8" = "4.15%"
    "D10" = "0.9456"
    "D11" = "0.1301"
    "E11" = "2.72%"
    "E12" = "3.73%"
    "D13" = "8.922"
    "E13" = "37.12%"
    "D14" = "0.09323"
    "E14" = "4.41%"
    "D15" = "0.03492"
    "E15" = "4.33%"
    "D16" = "0.09621"
    "E16" = "1.08%"
    "D17" = "0.001315"
    "E17" = "-5.46%"
    "D18" = "0.006286"
    "E18" = "2.62%"
    "D19" = "3.359"
    "E19" = "-0.95%"
    "E20" = "2.04%"
    "E21" = "9.39%"
    "D22" = "0.2413"
    "E22" = "5.11%"
    "E23" = "2.05%"
    "D24" = "0.001262"
    "E24" = "5.68%"
    "D25" = "0.004358"
    "E25" = "2.90%"
    "E26" = "-17.33%"
    "D27" = "0.0003997"
    "E27" = "1.11%"
    "E39" = "6.73%"
    "D40" = "0.05312"
    "E40" = "3.51%"
    "D41" = "0.007481"
    "E41" = "0.24%"
    "E42" = "3.35%"
    "D43" = "0.008662"
    "E43" = "2.56%"
    "D44" = "0.002052"
    "E44" = "3.36%"
    "D45" = "0.01047"
    "E45" = "32.48%"
    "D46" = "0.00006900"
    "E46" = "9.29%"
    "E47" = "0.86%"
    "D48" = "0.003507"
    "E48" = "23.22%"
    "D49" = "0.001703"
    "E49" = "1.50%"
    "D50" = "0.00002104"
    "E50" = "0.86%"
    "D51" = "0.0002004"
    "E51" = "0.86%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
